# Update countries & provincias Spain
# COVID-19 dataset refresh: new snapshot timestamp + updated counts for several
# countries (includes four countries that overtook a sort-neighbour in total cases).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" timestamp banner, row 1
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 14:07"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 4863077
$ws.Range("C4").Value = 903
$ws.Range("D4").Value = 2448305
$ws.Range("E4").Value = 2255797
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 46
$ws.Range("H4").Value = 158975

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1864561
$ws.Range("C6").Value = 9230
$ws.Range("D6").Value = 1235841
$ws.Range("E6").Value = 589663
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 86
$ws.Range("H6").Value = 39057

# Row 21: Alemania
$ws.Range("A21").Value = "Alemania"
$ws.Range("B21").Value = 212331
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 194700
$ws.Range("E21").Value = 8399
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 9232

# Row 40: Kuwait
$ws.Range("A40").Value = "Kuwait"
$ws.Range("B40").Value = 68774
$ws.Range("C40").Value = 475
$ws.Range("D40").Value = 60326
$ws.Range("E40").Value = 7983
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 465

# Row 41: Panama
$ws.Range("A41").Value = "Panama"
$ws.Range("B41").Value = 68456
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 42093
$ws.Range("E41").Value = 24866
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 1497

# Row 68: Nepal
$ws.Range("A68").Value = "Nepal"
$ws.Range("B68").Value = 21009
$ws.Range("C68").Value = 259
$ws.Range("D68").Value = 15026
$ws.Range("E68").Value = 5925
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 58

# Row 69: Venezuela
$ws.Range("A69").Value = "Venezuela"
$ws.Range("B69").Value = 20754
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 11622
$ws.Range("E69").Value = 8952
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 180

# Row 78: Dinamarca
$ws.Range("A78").Value = "Dinamarca"
$ws.Range("B78").Value = 14073
$ws.Range("C78").Value = 77
$ws.Range("D78").Value = 12715
$ws.Range("E78").Value = 742
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 616

# Row 93: Finlandia
$ws.Range("A93").Value = "Finlandia"
$ws.Range("B93").Value = 7483
$ws.Range("C93").Value = 17
$ws.Range("D93").Value = 6950
$ws.Range("E93").Value = 202
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 331

# Row 126: Eslovenia
$ws.Range("A126").Value = "Eslovenia"
$ws.Range("B126").Value = 2190
$ws.Range("C126").Value = 9
$ws.Range("D126").Value = 1854
$ws.Range("E126").Value = 213
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 123

# Row 132: Islandia
$ws.Range("A132").Value = "Islandia"
$ws.Range("B132").Value = 1918
$ws.Range("C132").Value = 3
$ws.Range("D132").Value = 1825
$ws.Range("E132").Value = 83
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 10

# Row 147: Burkina Faso
$ws.Range("A147").Value = "Burkina Faso"
$ws.Range("B147").Value = 1153
$ws.Range("C147").Value = 3
$ws.Range("D147").Value = 947
$ws.Range("E147").Value = 153
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 53

# Row 148: Niger
$ws.Range("A148").Value = "Niger"
$ws.Range("B148").Value = 1152
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 1032
$ws.Range("E148").Value = 51
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 69

# Row 161: Vietnam
$ws.Range("A161").Value = "Vietnam"
$ws.Range("B161").Value = 670
$ws.Range("C161").Value = 18
$ws.Range("D161").Value = 378
$ws.Range("E161").Value = 284
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 2
$ws.Range("H161").Value = 8

# Row 162: Reunion
$ws.Range("A162").Value = "Reunion"
$ws.Range("B162").Value = 667
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 592
$ws.Range("E162").Value = 71
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 4

